$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.091.04"
$ws.Range("E2").Value = "  +4.44%  "

# Row 3
$ws.Range("D3").Value = "3.488.71"
$ws.Range("E3").Value = "  +1.76%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "417.38"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.03"
$ws.Range("E6").Value = "  +1.76%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.653"
$ws.Range("E7").Value = "  +4.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.774"
$ws.Range("E9").Value = "  +6.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +13.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.94"
$ws.Range("E11").Value = "  +0.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000263"
$ws.Range("E12").Value = "  +18.10%  "

# Row 13
$ws.Range("E13").Value = "  +7.63%  "

# Row 14
$ws.Range("D14").Value = "4.038.08"
$ws.Range("E14").Value = "  +1.68%  "

# Row 15
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.31"
$ws.Range("E16").Value = "  -0.78%  "

# Row 17
$ws.Range("D17").Value = "3.513.33"
$ws.Range("E17").Value = "  +2.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.88"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("E19").Value = "  +2.18%  "

# Row 20
$ws.Range("D20").Value = "64.913.96"
$ws.Range("E20").Value = "  +4.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.48"
$ws.Range("E21").Value = "  -7.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.45"
$ws.Range("E22").Value = "  -2.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  -1.24%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.13"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.36"
$ws.Range("E25").Value = "  +1.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.89"
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.84"
$ws.Range("E27").Value = "  +1.03%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.81"
$ws.Range("E28").Value = "  +6.98%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.49"
$ws.Range("E29").Value = "  +5.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  -4.61%  "

# Row 31
$ws.Range("E31").Value = "  +4.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.161"
$ws.Range("E32").Value = "  -2.32%  "

# Row 33
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.34"
$ws.Range("E34").Value = "  -3.53%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "56.85"
$ws.Range("E35").Value = "  -2.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0503"
$ws.Range("E36").Value = "  +3.26%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0720"
$ws.Range("E37").Value = "  +31.76%  "

# Row 38
$ws.Range("E38").Value = "  +8.99%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").Value = "  -0.33%  "

# Row 40
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.78"
$ws.Range("E40").Value = "  +3.81%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.01"
$ws.Range("E41").Value = "  -0.53%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.47"
$ws.Range("E42").Value = "  +3.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.50"
$ws.Range("E43").Value = "  +1.70%  "

# Row 44
$ws.Range("E44").Value = "  -1.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.310"
$ws.Range("E45").Value = "  -4.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.99"
$ws.Range("E46").Value = "  -3.20%  "

# Row 47
$ws.Range("E47").Value = "  -2.85%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.147"
$ws.Range("E48").Value = "  +4.27%  "

# Row 49
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.66"
$ws.Range("E49").Value = "  -4.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("E50").Value = "  +9.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.37"
$ws.Range("E51").Value = "  -4.53%  "

